$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = [double]"8500"
$ws.Range("E2").Value = [double]"2.808192539087345E-37"
$ws.Range("F2").Value = [double]"85"
$ws.Range("H2").Value = [double]"24.26623452497573"
$ws.Range("I2").Value = [double]"19.11774104457871"
$ws.Range("D3").Value = [double]"3612.5"
$ws.Range("F3").Value = [double]"85"
$ws.Range("G3").Value = [double]"85"
$ws.Range("H3").Value = [double]"24.26623452497573"
$ws.Range("I3").Value = [double]"19.11774104457871"
$ws.Range("J3").Value = [double]"24.26623452497573"
$ws.Range("K3").Value = [double]"19.11774104457871"
$ws.Range("D4").Value = [double]"3512"
$ws.Range("E4").Value = [double]"0.5709916131736006"
$ws.Range("F4").Value = [double]"85"
$ws.Range("G4").Value = [double]"87"
$ws.Range("H4").Value = [double]"24.26623452497573"
$ws.Range("I4").Value = [double]"19.11774104457871"
$ws.Range("J4").Value = [double]"25.48038375748802"
$ws.Range("K4").Value = [double]"19.01076385234177"
$ws.Range("D5").Value = [double]"3144"
$ws.Range("E5").Value = [double]"0.1147689868577877"
$ws.Range("F5").Value = [double]"85"
$ws.Range("G5").Value = [double]"86"
$ws.Range("H5").Value = [double]"24.26623452497573"
$ws.Range("I5").Value = [double]"19.11774104457871"
$ws.Range("J5").Value = [double]"28.69016912565649"
$ws.Range("K5").Value = [double]"19.91412977917751"
$ws.Range("D6").Value = [double]"3675"
$ws.Range("E6").Value = [double]"0.4911218548390623"
$ws.Range("F6").Value = [double]"85"
$ws.Range("G6").Value = [double]"92"
$ws.Range("H6").Value = [double]"24.26623452497573"
$ws.Range("I6").Value = [double]"19.11774104457871"
$ws.Range("J6").Value = [double]"25.96458601641032"
$ws.Range("K6").Value = [double]"18.91157182213508"
$ws.Range("D7").Value = [double]"3556"
$ws.Range("E7").Value = [double]"0.4962400474312094"
$ws.Range("F7").Value = [double]"85"
$ws.Range("G7").Value = [double]"89"
$ws.Range("H7").Value = [double]"24.26623452497573"
$ws.Range("I7").Value = [double]"19.11774104457871"
$ws.Range("J7").Value = [double]"26.08149565073357"
$ws.Range("K7").Value = [double]"19.47935492619516"
$ws.Range("D8").Value = [double]"9400"
$ws.Range("E8").Value = [double]"2.565629657834091E-38"
$ws.Range("F8").Value = [double]"94"
$ws.Range("G8").Value = [double]"100"
$ws.Range("H8").Value = [double]"17.60975073558149"
$ws.Range("I8").Value = [double]"18.63055146756277"
$ws.Range("D9").Value = [double]"4418"
$ws.Range("F9").Value = [double]"94"
$ws.Range("G9").Value = [double]"94"
$ws.Range("H9").Value = [double]"17.60975073558149"
$ws.Range("I9").Value = [double]"18.63055146756277"
$ws.Range("J9").Value = [double]"17.60975073558149"
$ws.Range("K9").Value = [double]"18.63055146756277"
$ws.Range("D10").Value = [double]"3914"
$ws.Range("E10").Value = [double]"0.3194750096994206"
$ws.Range("F10").Value = [double]"94"
$ws.Range("G10").Value = [double]"91"
$ws.Range("H10").Value = [double]"17.60975073558149"
$ws.Range("I10").Value = [double]"18.63055146756277"
$ws.Range("J10").Value = [double]"19.45412092989661"
$ws.Range("K10").Value = [double]"19.08886489338698"
$ws.Range("D11").Value = [double]"4281"
$ws.Range("E11").Value = [double]"0.9078337969652105"
$ws.Range("F11").Value = [double]"94"
$ws.Range("G11").Value = [double]"92"
$ws.Range("H11").Value = [double]"17.60975073558149"
$ws.Range("I11").Value = [double]"18.63055146756277"
$ws.Range("J11").Value = [double]"17.52929804045222"
$ws.Range("K11").Value = [double]"18.09275629877325"
$ws.Range("D12").Value = [double]"4333"
$ws.Range("E12").Value = [double]"0.6376471556753459"
$ws.Range("F12").Value = [double]"94"
$ws.Range("G12").Value = [double]"96"
$ws.Range("H12").Value = [double]"17.60975073558149"
$ws.Range("I12").Value = [double]"18.63055146756277"
$ws.Range("J12").Value = [double]"17.6699092160549"
$ws.Range("K12").Value = [double]"18.01273644578071"
$ws.Range("D13").Value = [double]"4183"
$ws.Range("E13").Value = [double]"0.3860615910561968"
$ws.Range("F13").Value = [double]"94"
$ws.Range("G13").Value = [double]"96"
$ws.Range("H13").Value = [double]"17.60975073558149"
$ws.Range("I13").Value = [double]"18.63055146756277"
$ws.Range("J13").Value = [double]"20.69082520324423"
$ws.Range("K13").Value = [double]"22.07105363496164"
$ws.Range("D14").Value = [double]"7800"
$ws.Range("E14").Value = [double]"1.625080043660306E-26"
$ws.Range("F14").Value = [double]"85"
$ws.Range("H14").Value = [double]"1.056303268098712"
$ws.Range("I14").Value = [double]"0.04842710843895684"
$ws.Range("D15").Value = [double]"3612.5"
$ws.Range("F15").Value = [double]"85"
$ws.Range("G15").Value = [double]"85"
$ws.Range("H15").Value = [double]"1.056303268098712"
$ws.Range("I15").Value = [double]"0.04842710843895684"
$ws.Range("J15").Value = [double]"1.056303268098712"
$ws.Range("K15").Value = [double]"0.04842710843895684"
$ws.Range("D16").Value = [double]"3862"
$ws.Range("E16").Value = [double]"0.6154734676114035"
$ws.Range("F16").Value = [double]"85"
$ws.Range("G16").Value = [double]"87"
$ws.Range("H16").Value = [double]"1.056303268098712"
$ws.Range("I16").Value = [double]"0.04842710843895684"
$ws.Range("J16").Value = [double]"1.049903903745549"
$ws.Range("K16").Value = [double]"0.0621955229566507"
$ws.Range("D17").Value = [double]"4061"
$ws.Range("E17").Value = [double]"0.2103028377616685"
$ws.Range("F17").Value = [double]"85"
$ws.Range("G17").Value = [double]"86"
$ws.Range("H17").Value = [double]"1.056303268098712"
$ws.Range("I17").Value = [double]"0.04842710843895684"
$ws.Range("J17").Value = [double]"1.041633458452455"
$ws.Range("K17").Value = [double]"0.06606306204679072"
$ws.Range("D18").Value = [double]"4439"
$ws.Range("E18").Value = [double]"0.1207216408379425"
$ws.Range("F18").Value = [double]"85"
$ws.Range("G18").Value = [double]"92"
$ws.Range("H18").Value = [double]"1.056303268098712"
$ws.Range("I18").Value = [double]"0.04842710843895684"
$ws.Range("J18").Value = [double]"1.0396313861683"
$ws.Range("K18").Value = [double]"0.06156258878931065"
$ws.Range("D19").Value = [double]"3780"
$ws.Range("E19").Value = [double]"0.9951956508102576"
$ws.Range("F19").Value = [double]"85"
$ws.Range("G19").Value = [double]"89"
$ws.Range("H19").Value = [double]"1.056303268098712"
$ws.Range("I19").Value = [double]"0.04842710843895684"
$ws.Range("J19").Value = [double]"1.05550427863274"
$ws.Range("K19").Value = [double]"0.05466385591240686"
$ws.Range("D20").Value = [double]"4900"
$ws.Range("E20").Value = [double]"0.5826912833642663"
$ws.Range("F20").Value = [double]"94"
$ws.Range("G20").Value = [double]"100"
$ws.Range("H20").Value = [double]"1.000951915356464"
$ws.Range("I20").Value = [double]"0.02940571294469115"
$ws.Range("D21").Value = [double]"4418"
$ws.Range("F21").Value = [double]"94"
$ws.Range("G21").Value = [double]"94"
$ws.Range("H21").Value = [double]"1.000951915356464"
$ws.Range("I21").Value = [double]"0.02940571294469115"
$ws.Range("J21").Value = [double]"1.000951915356464"
$ws.Range("K21").Value = [double]"0.02940571294469115"
$ws.Range("D22").Value = [double]"4196"
$ws.Range("E22").Value = [double]"0.8250317469662701"
$ws.Range("F22").Value = [double]"94"
$ws.Range("G22").Value = [double]"91"
$ws.Range("H22").Value = [double]"1.000951915356464"
$ws.Range("I22").Value = [double]"0.02940571294469115"
$ws.Range("J22").Value = [double]"0.9994274149047458"
$ws.Range("K22").Value = [double]"0.02540785342055299"
$ws.Range("D23").Value = [double]"4611"
$ws.Range("E23").Value = [double]"0.4351350131115416"
$ws.Range("F23").Value = [double]"94"
$ws.Range("G23").Value = [double]"92"
$ws.Range("H23").Value = [double]"1.000951915356464"
$ws.Range("I23").Value = [double]"0.02940571294469115"
$ws.Range("J23").Value = [double]"0.9977133013073197"
$ws.Range("K23").Value = [double]"0.03125245642577462"
$ws.Range("D24").Value = [double]"5541"
$ws.Range("E24").Value = [double]"0.006651523926938394"
$ws.Range("F24").Value = [double]"94"
$ws.Range("G24").Value = [double]"96"
$ws.Range("H24").Value = [double]"1.000951915356464"
$ws.Range("I24").Value = [double]"0.02940571294469115"
$ws.Range("J24").Value = [double]"0.9872557777411911"
$ws.Range("K24").Value = [double]"0.03648096747536625"
$ws.Range("D25").Value = [double]"5312"
$ws.Range("E25").Value = [double]"0.03489586151989286"
$ws.Range("F25").Value = [double]"94"
$ws.Range("G25").Value = [double]"96"
$ws.Range("H25").Value = [double]"1.000951915356464"
$ws.Range("I25").Value = [double]"0.02940571294469115"
$ws.Range("J25").Value = [double]"0.9909675244626969"
$ws.Range("K25").Value = [double]"0.02975630978394025"
$ws.Range("D26").Value = [double]"1600"
$ws.Range("E26").Value = [double]"1.798394094827535E-15"
$ws.Range("F26").Value = [double]"85"
$ws.Range("H26").Value = [double]"0.7310366159814883"
$ws.Range("I26").Value = [double]"0.4251270990576909"
$ws.Range("D27").Value = [double]"3612.5"
$ws.Range("F27").Value = [double]"85"
$ws.Range("G27").Value = [double]"85"
$ws.Range("H27").Value = [double]"0.7310366159814883"
$ws.Range("I27").Value = [double]"0.4251270990576909"
$ws.Range("J27").Value = [double]"0.7310366159814883"
$ws.Range("K27").Value = [double]"0.4251270990576909"
$ws.Range("D28").Value = [double]"3613"
$ws.Range("E28").Value = [double]"0.7969755525344451"
$ws.Range("F28").Value = [double]"85"
$ws.Range("G28").Value = [double]"87"
$ws.Range("H28").Value = [double]"0.7310366159814883"
$ws.Range("I28").Value = [double]"0.4251270990576909"
$ws.Range("J28").Value = [double]"0.7629962790402484"
$ws.Range("K28").Value = [double]"0.4643945638808321"
$ws.Range("D29").Value = [double]"3899"
$ws.Range("E29").Value = [double]"0.451896064241222"
$ws.Range("F29").Value = [double]"85"
$ws.Range("G29").Value = [double]"86"
$ws.Range("H29").Value = [double]"0.7310366159814883"
$ws.Range("I29").Value = [double]"0.4251270990576909"
$ws.Range("J29").Value = [double]"0.705603482581682"
$ws.Range("K29").Value = [double]"0.4574071535320576"
$ws.Range("D30").Value = [double]"3453"
$ws.Range("E30").Value = [double]"0.1801317750976587"
$ws.Range("F30").Value = [double]"85"
$ws.Range("G30").Value = [double]"92"
$ws.Range("H30").Value = [double]"0.7310366159814883"
$ws.Range("I30").Value = [double]"0.4251270990576909"
$ws.Range("J30").Value = [double]"0.871156715925275"
$ws.Range("K30").Value = [double]"0.5848175195155606"
$ws.Range("D31").Value = [double]"3722"
$ws.Range("E31").Value = [double]"0.8566487012274508"
$ws.Range("F31").Value = [double]"85"
$ws.Range("G31").Value = [double]"89"
$ws.Range("H31").Value = [double]"0.7310366159814883"
$ws.Range("I31").Value = [double]"0.4251270990576909"
$ws.Range("J31").Value = [double]"0.765177299786501"
$ws.Range("K31").Value = [double]"0.4753618451576108"
$ws.Range("D32").Value = [double]"1400"
$ws.Range("E32").Value = [double]"1.014362195804016E-19"
$ws.Range("F32").Value = [double]"94"
$ws.Range("G32").Value = [double]"100"
$ws.Range("H32").Value = [double]"0.7152160735701159"
$ws.Range("I32").Value = [double]"0.4169546801315465"
$ws.Range("D33").Value = [double]"4418"
$ws.Range("F33").Value = [double]"94"
$ws.Range("G33").Value = [double]"94"
$ws.Range("H33").Value = [double]"0.7152160735701159"
$ws.Range("I33").Value = [double]"0.4169546801315465"
$ws.Range("J33").Value = [double]"0.7152160735701159"
$ws.Range("K33").Value = [double]"0.4169546801315465"
$ws.Range("D34").Value = [double]"4117"
$ws.Range("E34").Value = [double]"0.6613598866143098"
$ws.Range("F34").Value = [double]"94"
$ws.Range("G34").Value = [double]"91"
$ws.Range("H34").Value = [double]"0.7152160735701159"
$ws.Range("I34").Value = [double]"0.4169546801315465"
$ws.Range("J34").Value = [double]"0.7361271814716577"
$ws.Range("K34").Value = [double]"0.4082408668892568"
$ws.Range("D35").Value = [double]"4494"
$ws.Range("E35").Value = [double]"0.6442798432015006"
$ws.Range("F35").Value = [double]"94"
$ws.Range("G35").Value = [double]"92"
$ws.Range("H35").Value = [double]"0.7152160735701159"
$ws.Range("I35").Value = [double]"0.4169546801315465"
$ws.Range("J35").Value = [double]"0.7213559464104189"
$ws.Range("K35").Value = [double]"0.4795628221458288"
$ws.Range("D36").Value = [double]"4490"
$ws.Range("E36").Value = [double]"0.9547602684096348"
$ws.Range("F36").Value = [double]"94"
$ws.Range("G36").Value = [double]"96"
$ws.Range("H36").Value = [double]"0.7152160735701159"
$ws.Range("I36").Value = [double]"0.4169546801315465"
$ws.Range("J36").Value = [double]"0.7353176517681007"
$ws.Range("K36").Value = [double]"0.4785321109872692"
$ws.Range("D37").Value = [double]"5113"
$ws.Range("E37").Value = [double]"0.1130839933888709"
$ws.Range("F37").Value = [double]"94"
$ws.Range("G37").Value = [double]"96"
$ws.Range("H37").Value = [double]"0.7152160735701159"
$ws.Range("I37").Value = [double]"0.4169546801315465"
$ws.Range("J37").Value = [double]"0.6424708716332671"
$ws.Range("K37").Value = [double]"0.3875645555753933"
